$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "ELC_HV" -> "ELEC_HV" in B10 and apply the new "grey text box"
# style (Normal 2 cell style + custom Number Format "Text" + light-grey
# fill) that the author introduced for these label cells.
$b10 = $ws.Range("B10")
$b10.Value = "ELEC_HV"
$b10.Style = "Normal 2"
$b10.NumberFormat = "\Te\x\t"
$b10.Interior.Color = 15921906

# Propagate the exact same formatting (same underlying style record) to the
# other cells that need it by copying B10's format instead of rebuilding it
# property-by-property on each range - keeps a single shared cell style
# instead of creating duplicate near-identical ones.
$b10.Copy()
$ws.Range("B11:B13").PasteSpecial(-4122)
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B11").Value = "ELEC_HV"
$ws.Range("B12").Value = "ELEC_HV"
$ws.Range("B13").Value = "ELEC_HV"
$ws.Range("C4").Value = "ELEC_HV"

# Update the saved view state: drop the old scroll/selection and leave the
# cursor on D18 (matches the workbook's last-saved UI state in the diff).
[void]$ws.Range("D18").Select()
